$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.901.31"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "1.772.20"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4493"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07470"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.097"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.058"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.224"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "1.772.88"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001062"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06429"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.831"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "27.924.84"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "1.975.85"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.201"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.104"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09175"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.576"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.636"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02299"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2098"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6356"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.974"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.188"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.392"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.926"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.741"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5885"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.963"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06924"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.20%  "

Write-Host "Updated cryptos list"